$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.04920213766858694
$ws.Cells.Item(3, 3).Value = 0.04369851331357211
$ws.Cells.Item(4, 3).Value = 0.04033301125021183
$ws.Cells.Item(5, 3).Value = 0.03896499393145803
$ws.Cells.Item(6, 3).Value = 0.03873804407757575
$ws.Cells.Item(7, 3).Value = 0.04031454771903498
$ws.Cells.Item(8, 3).Value = 0.04730163079591421
$ws.Cells.Item(9, 3).Value = 0.06111325427887948
$ws.Cells.Item(10, 3).Value = 0.07132998657696987
$ws.Cells.Item(11, 3).Value = 0.07599344777028705
$ws.Cells.Item(12, 3).Value = 0.0777616744467764
$ws.Cells.Item(13, 3).Value = 0.07738075400068567
$ws.Cells.Item(14, 3).Value = 0.07613887524323104
$ws.Cells.Item(15, 3).Value = 0.0753784849450625
$ws.Cells.Item(16, 3).Value = 0.07102553677135859
$ws.Cells.Item(17, 3).Value = 0.06835919592811024
$ws.Cells.Item(18, 3).Value = 0.0668270760894103
$ws.Cells.Item(19, 3).Value = 0.0663085823388343
$ws.Cells.Item(20, 3).Value = 0.06864287836469884
$ws.Cells.Item(21, 3).Value = 0.07650358339031982
$ws.Cells.Item(22, 3).Value = 0.08165426983026691
$ws.Cells.Item(23, 3).Value = 0.07890403802888102
$ws.Cells.Item(24, 3).Value = 0.06851462306069322
$ws.Cells.Item(25, 3).Value = 0.05736483345238241
$ws.Cells.Item(2, 4).Value = 0.2515540573494377
$ws.Cells.Item(3, 4).Value = 0.2436303406308582
$ws.Cells.Item(4, 4).Value = 0.238834987617409
$ws.Cells.Item(5, 4).Value = 0.2368985575996732
$ws.Cells.Item(6, 4).Value = 0.2365780893588152
$ws.Cells.Item(7, 4).Value = 0.2388088003121709
$ws.Cells.Item(8, 4).Value = 0.2488075530524867
$ws.Cells.Item(9, 4).Value = 0.2689636792527068
$ws.Cells.Item(10, 4).Value = 0.2841010920807889
$ws.Cells.Item(11, 4).Value = 0.2910577752062409
$ws.Cells.Item(12, 4).Value = 0.2937021147619987
$ws.Cells.Item(13, 4).Value = 0.2931321669041438
$ws.Cells.Item(14, 4).Value = 0.2912751272022547
$ws.Cells.Item(15, 4).Value = 0.2901389337116314
$ws.Cells.Item(16, 4).Value = 0.2836478636222921
$ws.Cells.Item(17, 4).Value = 0.2796837719735521
$ws.Cells.Item(18, 4).Value = 0.2774103856657177
$ws.Cells.Item(19, 4).Value = 0.2766418038275162
$ws.Cells.Item(20, 4).Value = 0.2801050687389903
$ws.Cells.Item(21, 4).Value = 0.2918203147749239
$ws.Cells.Item(22, 4).Value = 0.2995350894694297
$ws.Cells.Item(23, 4).Value = 0.2954122968492072
$ws.Cells.Item(24, 4).Value = 0.2799145829469154
$ws.Cells.Item(25, 4).Value = 0.2634527787076877
$ws.Cells.Item(2, 5).Value = 0.1816556173362613
$ws.Cells.Item(3, 5).Value = 0.1769017472624768
$ws.Cells.Item(4, 5).Value = 0.1740649030444494
$ws.Cells.Item(5, 5).Value = 0.1729295520455914
$ws.Cells.Item(6, 5).Value = 0.1727422788545496
$ws.Cells.Item(7, 5).Value = 0.1740495074748836
$ws.Cells.Item(8, 5).Value = 0.1799994833695422
$ws.Cells.Item(9, 5).Value = 0.1923168520597685
$ws.Cells.Item(10, 5).Value = 0.2017614573111004
$ws.Cells.Item(11, 5).Value = 0.20614371848378
$ws.Cells.Item(12, 5).Value = 0.2078154790563786
$ws.Cells.Item(13, 5).Value = 0.2074548896046267
$ws.Cells.Item(14, 5).Value = 0.206281009160584
$ws.Cells.Item(15, 5).Value = 0.2055635728683427
$ws.Cells.Item(16, 5).Value = 0.2014767894286393
$ws.Cells.Item(17, 5).Value = 0.1989916334604374
$ws.Cells.Item(18, 5).Value = 0.1975703224838838
$ws.Cells.Item(19, 5).Value = 0.1970904812542855
$ws.Cells.Item(20, 5).Value = 0.1992553461778144
$ws.Cells.Item(21, 5).Value = 0.2066254730893178
$ws.Cells.Item(22, 5).Value = 0.2115139076193913
$ws.Cells.Item(23, 5).Value = 0.2088983215386619
$ws.Cells.Item(24, 5).Value = 0.1991360984859654
$ws.Cells.Item(25, 5).Value = 0.188915240978659
$ws.Cells.Item(2, 6).Value = 1.014992828093639
$ws.Cells.Item(3, 6).Value = 1.022889721210234
$ws.Cells.Item(4, 6).Value = 1.028548644049252
$ws.Cells.Item(5, 6).Value = 1.031058068167823
$ws.Cells.Item(6, 6).Value = 1.031487028323561
$ws.Cells.Item(7, 6).Value = 1.028581664080136
$ws.Cells.Item(8, 6).Value = 1.017547300474604
$ws.Cells.Item(9, 6).Value = 1.002354557958782
$ws.Cells.Item(10, 6).Value = 0.9951463862798988
$ws.Cells.Item(11, 6).Value = 0.9927310814958119
$ws.Cells.Item(12, 6).Value = 0.9919410871926928
$ws.Cells.Item(13, 6).Value = 0.9921056769180936
$ws.Cells.Item(14, 6).Value = 0.9926635880002266
$ws.Cells.Item(15, 6).Value = 0.9930215676620833
$ws.Cells.Item(16, 6).Value = 0.9953216490201413
$ws.Cells.Item(17, 6).Value = 0.9969541998085916
$ws.Cells.Item(18, 6).Value = 0.99797447271294
$ws.Cells.Item(19, 6).Value = 0.9983338652727198
$ws.Cells.Item(20, 6).Value = 0.9967719974589073
$ws.Cells.Item(21, 6).Value = 0.9924963302041334
$ws.Cells.Item(22, 6).Value = 0.9904285652922198
$ws.Cells.Item(23, 6).Value = 0.9914655432663864
$ws.Cells.Item(24, 6).Value = 0.9968541166879987
$ws.Cells.Item(25, 6).Value = 1.005772009633056
$ws.Cells.Item(2, 7).Value = 0.478553821415332
$ws.Cells.Item(3, 7).Value = 0.4857757305474379
$ws.Cells.Item(4, 7).Value = 0.4907462072019442
$ws.Cells.Item(5, 7).Value = 0.4929060660251494
$ws.Cells.Item(6, 7).Value = 0.4932728091635781
$ws.Cells.Item(7, 7).Value = 0.4907747924788168
$ws.Cells.Item(8, 7).Value = 0.4809323594509749
$ws.Cells.Item(9, 7).Value = 0.4659075040396061
$ws.Cells.Item(10, 7).Value = 0.4575054740944893
$ws.Cells.Item(11, 7).Value = 0.4542621154464257
$ws.Cells.Item(12, 7).Value = 0.4531176860936483
$ws.Cells.Item(13, 7).Value = 0.4533604256881887
$ws.Cells.Item(14, 7).Value = 0.4541662808229319
$ws.Cells.Item(15, 7).Value = 0.4546708142314344
$ws.Cells.Item(16, 7).Value = 0.4577291271877186
$ws.Cells.Item(17, 7).Value = 0.4597539199145118
$ws.Cells.Item(18, 7).Value = 0.4609729716004338
$ws.Cells.Item(19, 7).Value = 0.4613950561385991
$ws.Cells.Item(20, 7).Value = 0.4595327385975381
$ws.Cells.Item(21, 7).Value = 0.453927304160203
$ws.Cells.Item(22, 7).Value = 0.4507522937177129
$ws.Cells.Item(23, 7).Value = 0.4524019859447321
$ws.Cells.Item(24, 7).Value = 0.4596325634777898
$ws.Cells.Item(25, 7).Value = 0.4695112406151623
$ws.Cells.Item(2, 8).Value = 0.6347061767780957
$ws.Cells.Item(3, 8).Value = 0.6432981877701778
$ws.Cells.Item(4, 8).Value = 0.6489934932226902
$ws.Cells.Item(5, 8).Value = 0.6514198081027587
$ws.Cells.Item(6, 8).Value = 0.6518290600032515
$ws.Cells.Item(7, 8).Value = 0.6490257886320592
$ws.Cells.Item(8, 8).Value = 0.6375815034855705
$ws.Cells.Item(9, 8).Value = 0.6184754942592932
$ws.Cells.Item(10, 8).Value = 0.6064795031239996
$ws.Cells.Item(11, 8).Value = 0.601467022168535
$ws.Cells.Item(12, 8).Value = 0.5996329913386091
$ws.Cells.Item(13, 8).Value = 0.6000251298573716
$ws.Cells.Item(14, 8).Value = 0.6013148500330772
$ws.Cells.Item(15, 8).Value = 0.6021131916456994
$ws.Cells.Item(16, 8).Value = 0.6068160398506208
$ws.Cells.Item(17, 8).Value = 0.6098150604197201
$ws.Cells.Item(18, 8).Value = 0.6115818518244893
$ws.Cells.Item(19, 8).Value = 0.6121872371955561
$ws.Cells.Item(20, 8).Value = 0.6094914787245784
$ws.Cells.Item(21, 8).Value = 0.6009342873763899
$ws.Cells.Item(22, 8).Value = 0.5957152900856357
$ws.Cells.Item(23, 8).Value = 0.5984665277907197
$ws.Cells.Item(24, 8).Value = 0.6096376372769257
$ws.Cells.Item(25, 8).Value = 0.6232862519151325
$ws.Cells.Item(2, 10).Value = 0.1769492012324463
$ws.Cells.Item(3, 10).Value = 0.173096387539708
$ws.Cells.Item(4, 10).Value = 0.1708381363278164
$ws.Cells.Item(5, 10).Value = 0.1699448747971388
$ws.Cells.Item(6, 10).Value = 0.1697981794894048
$ws.Cells.Item(7, 10).Value = 0.1708259802053718
$ws.Cells.Item(8, 10).Value = 0.1755984538696538
$ws.Cells.Item(9, 10).Value = 0.1858105799349943
$ws.Cells.Item(10, 10).Value = 0.1938363863392709
$ws.Cells.Item(11, 10).Value = 0.1976017472509426
$ws.Cells.Item(12, 10).Value = 0.1990440699373579
$ws.Cells.Item(13, 10).Value = 0.1987327075502634
$ws.Cells.Item(14, 10).Value = 0.1977200779966353
$ws.Cells.Item(15, 10).Value = 0.1971019574494477
$ws.Cells.Item(16, 10).Value = 0.1935926140183994
$ws.Cells.Item(17, 10).Value = 0.1914690514896336
$ws.Cells.Item(18, 10).Value = 0.1902583997465825
$ws.Cells.Item(19, 10).Value = 0.1898503423740436
$ws.Cells.Item(20, 10).Value = 0.1916939940790314
$ws.Cells.Item(21, 10).Value = 0.1980170649488997
$ws.Cells.Item(22, 10).Value = 0.202245518404979
$ws.Cells.Item(23, 10).Value = 0.1999799280289949
$ws.Cells.Item(24, 10).Value = 0.1915922657255749
$ws.Cells.Item(25, 10).Value = 0.1829562912968896
$ws.Cells.Item(2, 11).Value = 1.846878404809047
$ws.Cells.Item(3, 11).Value = 1.626112624306813
$ws.Cells.Item(4, 11).Value = 1.490070283032424
$ws.Cells.Item(5, 11).Value = 1.434512272977543
$ws.Cells.Item(6, 11).Value = 1.42527979616608
$ws.Cells.Item(7, 11).Value = 1.489321487566542
$ws.Cells.Item(8, 11).Value = 1.770862571761995
$ws.Cells.Item(9, 11).Value = 2.318923225006074
$ws.Cells.Item(10, 11).Value = 2.718969382682417
$ws.Cells.Item(11, 11).Value = 2.900364375999288
$ws.Cells.Item(12, 11).Value = 2.968966170496344
$ws.Cells.Item(13, 11).Value = 2.954195548131281
$ws.Cells.Item(14, 11).Value = 2.906010084368518
$ws.Cells.Item(15, 11).Value = 2.876483446524048
$ws.Cells.Item(16, 11).Value = 2.70710261813332
$ws.Cells.Item(17, 11).Value = 2.603039502393813
$ws.Cells.Item(18, 11).Value = 2.543130040035578
$ws.Cells.Item(19, 11).Value = 2.522836382803746
$ws.Cells.Item(20, 11).Value = 2.614122931566612
$ws.Cells.Item(21, 11).Value = 2.920165758013468
$ws.Cells.Item(22, 11).Value = 3.119664464093319
$ws.Cells.Item(23, 11).Value = 3.013236952555133
$ws.Cells.Item(24, 11).Value = 2.609112368956175
$ws.Cells.Item(25, 11).Value = 2.171106648573982
$ws.Cells.Item(2, 15).Value = 2.186604846098476
$ws.Cells.Item(3, 15).Value = 2.219581923208082
$ws.Cells.Item(4, 15).Value = 2.241834676972502
$ws.Cells.Item(5, 15).Value = 2.251405663761616
$ws.Cells.Item(6, 15).Value = 2.253025249437982
$ws.Cells.Item(7, 15).Value = 2.241961720509465
$ws.Cells.Item(8, 15).Value = 2.197558495901575
$ws.Cells.Item(9, 15).Value = 2.126447932236488
$ws.Cells.Item(10, 15).Value = 2.08401431141445
$ws.Cells.Item(11, 15).Value = 2.066857564206003
$ws.Cells.Item(12, 15).Value = 2.060670817690351
$ws.Cells.Item(13, 15).Value = 2.061989427198569
$ws.Cells.Item(14, 15).Value = 2.066342352723865
$ws.Cells.Item(15, 15).Value = 2.069049078530327
$ws.Cells.Item(16, 15).Value = 2.085178864793889
$ws.Cells.Item(17, 15).Value = 2.095624791895887
$ws.Cells.Item(18, 15).Value = 2.101834965601341
$ws.Cells.Item(19, 15).Value = 2.10397226430058
$ws.Cells.Item(20, 15).Value = 2.09449189307341
$ws.Cells.Item(21, 15).Value = 2.065055364106087
$ws.Cells.Item(22, 15).Value = 2.047625312742412
$ws.Cells.Item(23, 15).Value = 2.056762102451131
$ws.Cells.Item(24, 15).Value = 2.095003439196404
$ws.Cells.Item(25, 15).Value = 2.143967962108647
